$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 -> prospectus 4455355: Commentarios column (C) gets duplicated with the
# same JSON text already present in the Json column (B11).
$ws.Range("C11").Value = $ws.Range("B11").Value2

# Row 13 -> prospectus 4606694: add comment about reported rate on a revolving debt.
$ws.Range("C13").Value = 'Revolvente con tasa reportada ¿es correcto? Trae una tasa del "tasa_externa": 0.0124'

# Row 14 -> prospectus 4688988: add comment about revolving debt with zero balance.
$ws.Range("C14").Value = 'Revolvente con saldo cero'

# Row 15 -> prospectus 4699870: add comment about revolving loans with zero balance.
$ws.Range("C15").Value = 'Creditos Revolventes con saldo cero'
